# Roll the 5-year reporting window forward by one year:
#   - column D..H headers (row 8: period-end labels, row 9: publish dates)
#     each shift one column to the left, and column H receives the newest
#     (1401/12) period's data
#   - all financial statement rows (11-26) shift the same way, with new
#     H-column figures for the newly published period
#   - row 15 ("هزینه کاهش ارزش دریافتنی‌ها") only ever has a single real
#     figure among the five periods; that figure moves from column E to D

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: financial-period column headers -------------------------------
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish-date column headers ------------------------------------
$ws.Range("D9").Value = "1399-04-08 (11)"
$ws.Range("E9").Value = "1400-04-06 (11)"
$ws.Range("F9").Value = "1401-02-21 (10)"
$ws.Range("G9").Value = "1402-02-23 (10)"
$ws.Range("H9").Value = "1402-02-23 (2)"

# --- Row 11: فروش (Sales) ---------------------------------------------------
$ws.Range("D11").Value = 45585
$ws.Range("E11").Value = 40080
$ws.Range("F11").Value = 39864
$ws.Range("G11").Value = 104958
$ws.Range("H11").Value = 259744

# --- Row 12: بهای تمام شده کالای فروش رفته (COGS) --------------------------
$ws.Range("D12").Value = -30025
$ws.Range("E12").Value = -14673
$ws.Range("F12").Value = -14749
$ws.Range("G12").Value = -49402
$ws.Range("H12").Value = -156697

# --- Row 13: سود (زیان) ناخالص (Gross profit) -------------------------------
$ws.Range("D13").Value = 15560
$ws.Range("E13").Value = 25408
$ws.Range("F13").Value = 25115
$ws.Range("G13").Value = 55556
$ws.Range("H13").Value = 103047

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) --------------
$ws.Range("D14").Value = -4889
$ws.Range("E14").Value = -4075
$ws.Range("F14").Value = -3968
$ws.Range("G14").Value = -7915
$ws.Range("H14").Value = -27906

# --- Row 15: هزینه کاهش ارزش دریافتنی‌ها (only one real figure) -----------
$ws.Range("D15").Value = -785
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = "-"

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی ------------------------
$ws.Range("D16").Value = 530
$ws.Range("E16").Value = -310
$ws.Range("F16").Value = 435
$ws.Range("G16").Value = 824
$ws.Range("H16").Value = 572

# --- Row 17: سود (زیان) عملیاتی (Operating profit) --------------------------
$ws.Range("D17").Value = 10417
$ws.Range("E17").Value = 21023
$ws.Range("F17").Value = 21581
$ws.Range("G17").Value = 48465
$ws.Range("H17").Value = 75712

# --- Row 18: هزینه های مالی (Finance costs) ---------------------------------
$ws.Range("D18").Value = -3251
$ws.Range("E18").Value = -1010
$ws.Range("F18").Value = -1423
$ws.Range("G18").Value = -7404
$ws.Range("H18").Value = -14856

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی -----------------------
$ws.Range("D19").Value = 130
$ws.Range("E19").Value = 1037
$ws.Range("F19").Value = 784
$ws.Range("G19").Value = 2711
$ws.Range("H19").Value = 1143

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات --------------
$ws.Range("D20").Value = 7296
$ws.Range("E20").Value = 21049
$ws.Range("F20").Value = 20943
$ws.Range("G20").Value = 43772
$ws.Range("H20").Value = 61999

# --- Row 21: مالیات (Tax) ----------------------------------------------------
$ws.Range("D21").Value = -1703
$ws.Range("E21").Value = -3949
$ws.Range("F21").Value = -2648
$ws.Range("G21").Value = -5921
$ws.Range("H21").Value = -12220

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم ----------------------------
$ws.Range("D22").Value = 5593
$ws.Range("E22").Value = 17100
$ws.Range("F22").Value = 18294
$ws.Range("G22").Value = 37850
$ws.Range("H22").Value = 49780

# --- Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی (all "-", unchanged)

# --- Row 24: سود (زیان) خالص (Net profit) -----------------------------------
$ws.Range("D24").Value = 5593
$ws.Range("E24").Value = 17100
$ws.Range("F24").Value = 18294
$ws.Range("G24").Value = 37850
$ws.Range("H24").Value = 49780

# --- Row 25: سود هر سهم پس از کسر مالیات (all zero, unchanged)

# --- Row 26: سرمایه (Capital) ------------------------------------------------
$ws.Range("D26").Value = 9886
$ws.Range("E26").Value = 7795
$ws.Range("F26").Value = 4423
$ws.Range("G26").Value = 64051
$ws.Range("H26").Value = 47890

# --- Row 27: سود هر سهم بر اساس آخرین سرمایه (all zero, unchanged)
